$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.369278686223538
$ws.Range("C2").Value = 0.06750128626842411
$ws.Range("D2").Value = 0.02643641156463872
$ws.Range("E2").Value = 0.4159123285804327
$ws.Range("F2").Value = 0.7099491187442979
$ws.Range("K2").Value = 0.3603870442761661
$ws.Range("O2").Value = 2.425853481378994
$ws.Range("B3").Value = 0.3267623426023363
$ws.Range("C3").Value = 0.06179291441249291
$ws.Range("D3").Value = 0.02481443763277014
$ws.Range("E3").Value = 0.3628992592812068
$ws.Range("F3").Value = 0.7073304550515758
$ws.Range("K3").Value = 0.3148674128965752
$ws.Range("O3").Value = 2.431571157312789
$ws.Range("B4").Value = 0.3006692369605446
$ws.Range("C4").Value = 0.05826342446459876
$ws.Range("D4").Value = 0.02381009723378469
$ws.Range("E4").Value = 0.3304359795848484
$ws.Range("F4").Value = 0.7062327119439757
$ws.Range("K4").Value = 0.2868644905344411
$ws.Range("O4").Value = 2.436779632266536
$ws.Range("B5").Value = 0.2900395173022901
$ws.Range("C5").Value = 0.05681901421743873
$ws.Range("D5").Value = 0.02339872983109359
$ws.Range("E5").Value = 0.3172271342781698
$ws.Range("F5").Value = 0.7059133797578028
$ws.Range("K5").Value = 0.2754399761403761
$ws.Range("O5").Value = 2.439328064159156
$ws.Range("B6").Value = 0.2882746793737567
$ws.Range("C6").Value = 0.05657880328010378
$ws.Range("D6").Value = 0.02333029727888913
$ws.Range("E6").Value = 0.3150349875118934
$ws.Range("F6").Value = 0.7058680783579803
$ws.Range("K6").Value = 0.2735421652580783
$ws.Range("O6").Value = 2.439776925921819
$ws.Range("B7").Value = 0.3005258664247776
$ws.Range("C7").Value = 0.05824396932864317
$ws.Range("D7").Value = 0.02380455781059965
$ws.Range("E7").Value = 0.3302577608487383
$ws.Range("F7").Value = 0.7062278873820702
$ws.Range("K7").Value = 0.2867104680244665
$ws.Range("O7").Value = 2.436812278036115
$ws.Range("B8").Value = 0.3546167587350624
$ws.Range("C8").Value = 0.06553816001786572
$ws.Range("D8").Value = 0.02587892582896956
$ws.Range("E8").Value = 0.3976143566367938
$ws.Range("F8").Value = 0.7089401619904407
$ws.Range("K8").Value = 0.3447032443233127
$ws.Range("O8").Value = 2.427471956429031
$ws.Range("B9").Value = 0.4607732036895413
$ws.Range("C9").Value = 0.0796457980496541
$ws.Range("D9").Value = 0.02987854830368519
$ws.Range("E9").Value = 0.5304715718473005
$ws.Range("F9").Value = 0.7183200586844052
$ws.Range("K9").Value = 0.4579889468427893
$ws.Range("O9").Value = 2.422674350482254
$ws.Range("B10").Value = 0.538810212772205
$ws.Range("C10").Value = 0.08988985053767351
$ws.Range("D10").Value = 0.03277412581060446
$ws.Range("E10").Value = 0.6286731617453114
$ws.Range("F10").Value = 0.7277079534868562
$ws.Range("K10").Value = 0.5409453036400578
$ws.Range("O10").Value = 2.42746075108019
$ws.Range("B11").Value = 0.5743197598497147
$ws.Range("C11").Value = 0.09452372209453586
$ws.Range("D11").Value = 0.03408180682819761
$ws.Range("E11").Value = 0.6735028481930811
$ws.Range("F11").Value = 0.7325254296817718
$ws.Range("K11").Value = 0.5786236393819877
$ws.Range("O11").Value = 2.431458261238959
$ws.Range("B12").Value = 0.5877674948827689
$ws.Range("C12").Value = 0.09627464524143647
$ws.Range("D12").Value = 0.0345755930863163
$ws.Range("E12").Value = 0.6905033027698693
$ws.Range("F12").Value = 0.7344286564715787
$ws.Range("K12").Value = 0.5928827030609227
$ws.Range("O12").Value = 2.4332349316457
$ws.Range("B13").Value = 0.5848712448268714
$ws.Range("C13").Value = 0.09589772354335935
$ws.Range("D13").Value = 0.03446931036020828
$ws.Range("E13").Value = 0.686840835426068
$ws.Range("F13").Value = 0.7340152460770497
$ws.Range("K13").Value = 0.5898121613913645
$ws.Range("O13").Value = 2.432840582747644
$ws.Range("B14").Value = 0.5754260946169438
$ws.Range("C14").Value = 0.09466784854060961
$ws.Range("D14").Value = 0.03412245925307644
$ws.Range("E14").Value = 0.674900985891199
$ws.Range("F14").Value = 0.7326804251625418
$ws.Range("K14").Value = 0.5797969220361665
$ws.Range("O14").Value = 2.43159915339146
$ws.Range("B15").Value = 0.5696407906200704
$ws.Range("C15").Value = 0.09391401501756036
$ws.Range("D15").Value = 0.03390981908308532
$ws.Range("E15").Value = 0.6675907181106737
$ws.Range("F15").Value = 0.7318731001349335
$ws.Range("K15").Value = 0.5736611264520377
$ws.Range("O15").Value = 2.430873014034944
$ws.Range("B16").Value = 0.5364897484143398
$ws.Range("C16").Value = 0.08958648412675529
$ws.Range("D16").Value = 0.03268847110424389
$ws.Range("E16").Value = 0.6257467554810319
$ws.Range("F16").Value = 0.7274041513734488
$ws.Range("K16").Value = 0.5384817130796762
$ws.Range("O16").Value = 2.427236229956179
$ws.Range("B17").Value = 0.5161549529773879
$ws.Range("C17").Value = 0.08692493365489895
$ws.Range("D17").Value = 0.03193674841013205
$ws.Range("E17").Value = 0.6001184060109779
$ws.Range("F17").Value = 0.7248028866342224
$ws.Range("K17").Value = 0.5168848815781359
$ws.Range("O17").Value = 2.425472192730808
$ws.Range("B18").Value = 0.5044598935531042
$ws.Range("C18").Value = 0.08539161750429969
$ws.Range("D18").Value = 0.03150348220748356
$ws.Range("E18").Value = 0.5853923531770278
$ws.Range("F18").Value = 0.7233581637360658
$ws.Range("K18").Value = 0.5044574342174997
$ws.Range("O18").Value = 2.424628816657702
$ws.Range("B19").Value = 0.500500327919724
$ws.Range("C19").Value = 0.08487204175173702
$ws.Range("D19").Value = 0.0313566329935
$ws.Range("E19").Value = 0.5804088339549196
$ws.Range("F19").Value = 0.7228778334412738
$ws.Range("K19").Value = 0.500248779034365
$ws.Range("O19").Value = 2.424372641221453
$ws.Range("B20").Value = 0.5183195285641773
$ws.Range("C20").Value = 0.0872085155361475
$ws.Range("D20").Value = 0.03201686346714894
$ws.Range("E20").Value = 0.6028450535715706
$ws.Range("F20").Value = 0.7250744679570147
$ws.Range("K20").Value = 0.5191844757123079
$ws.Range("O20").Value = 2.425642245181734
$ws.Range("B21").Value = 0.5782003392407375
$ws.Range("C21").Value = 0.09502919684966571
$ws.Range("D21").Value = 0.03422437620519503
$ws.Range("E21").Value = 0.6784073317457597
$ws.Range("F21").Value = 0.7330703491592914
$ws.Range("K21").Value = 0.5827388856736775
$ws.Range("O21").Value = 2.431956646954234
$ws.Range("B22").Value = 0.6173417845366771
$ws.Range("C22").Value = 0.1001181529368154
$ws.Range("D22").Value = 0.03565891356843309
$ws.Range("E22").Value = 0.727934868665173
$ws.Range("F22").Value = 0.7387564645732709
$ws.Range("K22").Value = 0.6242232693227834
$ws.Range("O22").Value = 2.437616399496534
$ws.Range("B23").Value = 0.5964508569768725
$ws.Range("C23").Value = 0.09740414103167438
$ws.Range("D23").Value = 0.03489403523912671
$ws.Range("E23").Value = 0.7014873838633662
$ws.Range("F23").Value = 0.7356794526076555
$ws.Range("K23").Value = 0.602087180787521
$ws.Range("O23").Value = 2.434455034506044
$ws.Range("B24").Value = 0.5173409373524009
$ws.Range("C24").Value = 0.08708031799307037
$ws.Range("D24").Value = 0.031980646845156
$ws.Range("E24").Value = 0.6016123113792702
$ws.Range("F24").Value = 0.7249515278863043
$ws.Range("K24").Value = 0.5181448638791153
$ws.Range("O24").Value = 2.42556483252423
$ws.Range("B25").Value = 0.43204684069525
$ws.Range("C25").Value = 0.07585044468351043
$ws.Range("D25").Value = 0.028803988992955
$ws.Range("E25").Value = 0.4944351383057182
$ws.Range("F25").Value = 0.7153456369876636
$ws.Range("K25").Value = 0.4273898620324417
$ws.Range("O25").Value = 2.439328064159156
